$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 9099.1
$ws.Range("H32").Value = 25012.908
$ws.Range("I32").Value = 42166.332
$ws.Range("J32").Value = 18580.375
$ws.Range("K32").Value = 42166.332
$ws.Range("L32").Value = 18580.375
$ws.Range("M32").Value = -41840.332
$ws.Range("N32").Value = -19232.375
$ws.Range("H33").Value = 758.35297
$ws.Range("I33").Value = 649.6667
$ws.Range("K33").Value = 649.6667
$ws.Range("M33").Value = -420.6667
$ws.Range("H74").Value = 8411304
$ws.Range("I74").Value = 14291967
$ws.Range("K74").Value = 14291967
$ws.Range("M74").Value = -14291031
$ws.Range("H77").Value = 8411304
$ws.Range("I77").Value = 14291967
$ws.Range("K77").Value = 71459835
$ws.Range("M77").Value = -71455155
$ws.Range("H100").Value = 3881
$ws.Range("I100").Value = 3428.625
$ws.Range("K100").Value = 3428.625
$ws.Range("M100").Value = -2887.625
$ws.Range("H101").Value = 1880.8334
$ws.Range("I101").Value = 2537.75
$ws.Range("J101").Value = 567
$ws.Range("K101").Value = 7613.25
$ws.Range("L101").Value = 1701
$ws.Range("M101").Value = -5991.25
$ws.Range("N101").Value = -4945
$ws.Range("H132").Value = 288152.25
$ws.Range("I132").Value = 320032.78
$ws.Range("K132").Value = 960098.3400000001
$ws.Range("M132").Value = -957568.3400000001
$ws.Range("H137").Value = 5703.706
$ws.Range("I137").Value = 5656.273
$ws.Range("J137").Value = 5790.6665
$ws.Range("K137").Value = 16968.819
$ws.Range("L137").Value = 17371.9995
$ws.Range("M137").Value = -14418.819
$ws.Range("N137").Value = -22471.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 199.61539
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H10").Value = 50000000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H11").Value = 12502999
$ws.Range("I11").Value = 16669332
$ws.Range("K11").Value = 16669332
$ws.Range("M11").Value = -16669188
$ws.Range("H32").Value = 1373579.9
$ws.Range("I32").Value = 1309.7922
$ws.Range("K32").Value = 1309.7922
$ws.Range("M32").Value = -1022.7922
$ws.Range("H61").Value = 4852.123
$ws.Range("I61").Value = 5381.5
$ws.Range("K61").Value = 5381.5
$ws.Range("M61").Value = -5169.5
$ws.Range("H63").Value = 4805.55
$ws.Range("I63").Value = 3601.0833
$ws.Range("K63").Value = 3601.0833
$ws.Range("M63").Value = -2915.0833
$ws.Range("H66").Value = 4805.55
$ws.Range("I66").Value = 3601.0833
$ws.Range("K66").Value = 18005.4165
$ws.Range("M66").Value = -14573.4165
$ws.Range("H68").Value = 13400
$ws.Range("I68").Value = 11000
$ws.Range("J68").Value = 15800
$ws.Range("K68").Value = 11000
$ws.Range("L68").Value = 15800
$ws.Range("M68").Value = -10189
$ws.Range("N68").Value = -17422
$ws.Range("H71").Value = 13400
$ws.Range("I71").Value = 11000
$ws.Range("J71").Value = 15800
$ws.Range("K71").Value = 33000
$ws.Range("L71").Value = 47400
$ws.Range("M71").Value = -28944
$ws.Range("N71").Value = -55512
$ws.Range("H93").Value = 8500
$ws.Range("I93").Value = 8500
$ws.Range("K93").Value = 8500
$ws.Range("M93").Value = -6004
$ws.Range("H132").Value = 617158.6
$ws.Range("I132").Value = 715649.1
$ws.Range("J132").Value = 105008.3
$ws.Range("K132").Value = 2146947.3
$ws.Range("L132").Value = 315024.9
$ws.Range("M132").Value = -2144417.3
$ws.Range("N132").Value = -320084.9
$ws.Range("H136").Value = 4852.123
$ws.Range("I136").Value = 5381.5
$ws.Range("K136").Value = 16144.5
$ws.Range("M136").Value = -13594.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 199.61539
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H95").Value = 26096.334
$ws.Range("J95").Value = 26096.334
$ws.Range("L95").Value = 26096.334
$ws.Range("N95").Value = -31588.334
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H105").Value = 2064.4375
$ws.Range("I105").Value = 2362.5833
$ws.Range("J105").Value = 1170
$ws.Range("K105").Value = 2362.5833
$ws.Range("L105").Value = 1170
$ws.Range("M105").Value = -615.5832999999998
$ws.Range("N105").Value = -4664
$ws.Range("H130").Value = 75000
$ws.Range("J130").Value = 75000
$ws.Range("L130").Value = 75000
$ws.Range("N130").Value = -85040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2001.3334
$ws.Range("J2").Value = 2250
$ws.Range("L2").Value = 2250
$ws.Range("N2").Value = -2476
$ws.Range("H6").Value = 25200
$ws.Range("I6").Value = 30090
$ws.Range("J6").Value = 750
$ws.Range("K6").Value = 30090
$ws.Range("L6").Value = 750
$ws.Range("M6").Value = -29977
$ws.Range("N6").Value = -976
$ws.Range("H22").Value = 777529.75
$ws.Range("I22").Value = 1116424
$ws.Range("J22").Value = 2914.4285
$ws.Range("K22").Value = 1116424
$ws.Range("L22").Value = 2914.4285
$ws.Range("M22").Value = -1116074
$ws.Range("N22").Value = -3614.4285
$ws.Range("H107").Value = 728.53845
$ws.Range("I107").Value = 616
$ws.Range("J107").Value = 825
$ws.Range("K107").Value = 616
$ws.Range("L107").Value = 825
$ws.Range("M107").Value = 1304
$ws.Range("N107").Value = -4665
$ws.Range("H122").Value = 7365.231
$ws.Range("I122").Value = 2081.5557
$ws.Range("J122").Value = 19253.5
$ws.Range("K122").Value = 6244.6671
$ws.Range("L122").Value = 57760.5
$ws.Range("M122").Value = -3794.6671
$ws.Range("N122").Value = -62660.5
$ws.Range("H132").Value = 7742.086
$ws.Range("I132").Value = 3514.8276
$ws.Range("K132").Value = 10544.4828
$ws.Range("M132").Value = -8014.4828

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 82.28
$ws.Range("I2").Value = 36.666668
$ws.Range("J2").Value = 124.38461
$ws.Range("K2").Value = 220.000008
$ws.Range("L2").Value = 746.3076599999999
$ws.Range("M2").Value = -107.000008
$ws.Range("N2").Value = -972.3076599999999
$ws.Range("H5").Value = 1076.9131
$ws.Range("I5").Value = 657.5
$ws.Range("K5").Value = 1972.5
$ws.Range("M5").Value = -1860.5
$ws.Range("H22").Value = 1254.3334
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1331
$ws.Range("H27").Value = 1254.3334
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 1500
$ws.Range("M27").Value = -1398
$ws.Range("H98").Value = 1040.5
$ws.Range("J98").Value = 1064.7778
$ws.Range("L98").Value = 3194.3334
$ws.Range("N98").Value = -6190.3334
$ws.Range("H114").Value = 1806.8572
$ws.Range("I114").Value = 984.6667
$ws.Range("K114").Value = 2954.0001
$ws.Range("M114").Value = 299.9998999999998
$ws.Range("H115").Value = 2004.6666
$ws.Range("I115").Value = 1405.6
$ws.Range("K115").Value = 4216.799999999999
$ws.Range("M115").Value = -3041.799999999999
$ws.Range("H130").Value = 10832.667
$ws.Range("I130").Value = 4999.6665
$ws.Range("J130").Value = 16665.666
$ws.Range("K130").Value = 14998.9995
$ws.Range("L130").Value = 49996.99800000001
$ws.Range("M130").Value = -9978.999500000002
$ws.Range("N130").Value = -60036.99800000001
$ws.Range("H135").Value = 1076.9131
$ws.Range("I135").Value = 657.5
$ws.Range("K135").Value = 5917.5
$ws.Range("M135").Value = -3382.5
$ws.Range("H140").Value = 32610982
$ws.Range("I140").Value = 39475296
$ws.Range("J140").Value = 5496.25
$ws.Range("K140").Value = 118425888
$ws.Range("L140").Value = 16488.75
$ws.Range("M140").Value = -118420708
$ws.Range("N140").Value = -26848.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -360
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H41").Value = 14850.429
$ws.Range("I41").Value = 16483
$ws.Range("J41").Value = 5055
$ws.Range("K41").Value = 16483
$ws.Range("L41").Value = 5055
$ws.Range("M41").Value = -16128
$ws.Range("N41").Value = -5765
$ws.Range("H97").Value = 2226.375
$ws.Range("I97").Value = 1846.3889
$ws.Range("J97").Value = 3366.3333
$ws.Range("K97").Value = 1846.3889
$ws.Range("L97").Value = 3366.3333
$ws.Range("M97").Value = -1350.3889
$ws.Range("N97").Value = -4358.3333
$ws.Range("H102").Value = 942675.9
$ws.Range("I102").Value = 1611048.8
$ws.Range("K102").Value = 1611048.8
$ws.Range("M102").Value = -1609426.8
$ws.Range("H132").Value = 4969.4253
$ws.Range("I132").Value = 5058.7295
$ws.Range("J132").Value = 4639
$ws.Range("K132").Value = 15176.1885
$ws.Range("L132").Value = 13917
$ws.Range("M132").Value = -12646.1885
$ws.Range("N132").Value = -18977

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 25000004
$ws.Range("J10").Value = 8
$ws.Range("L10").Value = 8
$ws.Range("N10").Value = -288
$ws.Range("H18").Value = 30000
$ws.Range("J18").Value = 30000
$ws.Range("L18").Value = 30000
$ws.Range("N18").Value = -30344
$ws.Range("H22").Value = 810
$ws.Range("I22").Value = 598.4
$ws.Range("J22").Value = 1162.6666
$ws.Range("K22").Value = 598.4
$ws.Range("L22").Value = 1162.6666
$ws.Range("M22").Value = -303.4
$ws.Range("N22").Value = -1752.6666
$ws.Range("H27").Value = 810
$ws.Range("I27").Value = 598.4
$ws.Range("J27").Value = 1162.6666
$ws.Range("K27").Value = 598.4
$ws.Range("L27").Value = 1162.6666
$ws.Range("M27").Value = -491.4
$ws.Range("N27").Value = -1376.6666
$ws.Range("H61").Value = 4870.75
$ws.Range("I61").Value = 3908.8845
$ws.Range("J61").Value = 17375
$ws.Range("K61").Value = 3908.8845
$ws.Range("L61").Value = 17375
$ws.Range("M61").Value = -3706.8845
$ws.Range("N61").Value = -17779
$ws.Range("H68").Value = 4493.7646
$ws.Range("I68").Value = 2533
$ws.Range("K68").Value = 2533
$ws.Range("M68").Value = -1784
$ws.Range("H71").Value = 4493.7646
$ws.Range("I71").Value = 2533
$ws.Range("K71").Value = 12665
$ws.Range("M71").Value = -8921
$ws.Range("H88").Value = 18900
$ws.Range("J88").Value = 18900
$ws.Range("L88").Value = 18900
$ws.Range("N88").Value = -19756
$ws.Range("H91").Value = 18900
$ws.Range("J91").Value = 18900
$ws.Range("L91").Value = 18900
$ws.Range("N91").Value = -21864
$ws.Range("H93").Value = 2711.4707
$ws.Range("I93").Value = 2926.2727
$ws.Range("J93").Value = 2317.6667
$ws.Range("K93").Value = 2926.2727
$ws.Range("L93").Value = 2317.6667
$ws.Range("M93").Value = -1678.2727
$ws.Range("N93").Value = -4813.6667
$ws.Range("H100").Value = 3577.111
$ws.Range("J100").Value = 2313.7144
$ws.Range("L100").Value = 2313.7144
$ws.Range("N100").Value = -3395.7144
$ws.Range("H113").Value = 4870.75
$ws.Range("I113").Value = 3908.8845
$ws.Range("J113").Value = 17375
$ws.Range("K113").Value = 3908.8845
$ws.Range("L113").Value = 17375
$ws.Range("M113").Value = -1738.8845
$ws.Range("N113").Value = -21715
$ws.Range("H132").Value = 2401.4614
$ws.Range("I132").Value = 1910.7894
$ws.Range("K132").Value = 5732.3682
$ws.Range("M132").Value = -3202.3682

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2500624.8
$ws.Range("I3").Value = 3333833
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 3333833
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -3333719
$ws.Range("N3").Value = -1228
$ws.Range("H9").Value = 5005000
$ws.Range("I9").Value = 5005000
$ws.Range("K9").Value = 5005000
$ws.Range("M9").Value = -5004860
$ws.Range("H51").Value = 17458.8
$ws.Range("J51").Value = 8612
$ws.Range("L51").Value = 8612
$ws.Range("N51").Value = -9632
$ws.Range("H81").Value = 1707.6111
$ws.Range("J81").Value = 2999.8
$ws.Range("L81").Value = 5999.6
$ws.Range("N81").Value = -8121.6
$ws.Range("H84").Value = 1707.6111
$ws.Range("J84").Value = 2999.8
$ws.Range("L84").Value = 29998
$ws.Range("N84").Value = -40606
$ws.Range("H132").Value = 8585.200000000001
$ws.Range("I132").Value = 5597.375
$ws.Range("K132").Value = 16792.125
$ws.Range("M132").Value = -14262.125
